$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A2:A4").EntireRow.Delete()
$ws.Range("A2:XFD4").Select() | Out-Null
